# Update "想去人数" (F column) counts across sheets, per commit
# "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 9804
$ws1.Range("F3").Value  = 416
$ws1.Range("F5").Value  = 18
$ws1.Range("F6").Value  = 275
$ws1.Range("F13").Value = 3036
$ws1.Range("F14").Value = 2297
$ws1.Range("F16").Value = 1992
$ws1.Range("F21").Value = 324
$ws1.Range("F22").Value = 30
$ws1.Range("F24").Value = 225
$ws1.Range("F28").Value = 331
$ws1.Range("F29").Value = 540
$ws1.Range("F31").Value = 181
$ws1.Range("F33").Value = 225
$ws1.Range("F34").Value = 1562
$ws1.Range("F35").Value = 74
$ws1.Range("F36").Value = 376
$ws1.Range("F38").Value = 404
$ws1.Range("F39").Value = 836
$ws1.Range("F41").Value = 327

# --- Sheet "演出" (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F9").Value = 6

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 9804
$ws4.Range("F3").Value  = 416
$ws4.Range("F6").Value  = 18
$ws4.Range("F8").Value  = 275
$ws4.Range("F15").Value = 3036
$ws4.Range("F16").Value = 2297
$ws4.Range("F18").Value = 1992
$ws4.Range("F23").Value = 324
$ws4.Range("F24").Value = 30
$ws4.Range("F26").Value = 225
$ws4.Range("F30").Value = 331
$ws4.Range("F31").Value = 540
$ws4.Range("F36").Value = 181
$ws4.Range("F39").Value = 225
$ws4.Range("F40").Value = 1562
$ws4.Range("F41").Value = 74
$ws4.Range("F43").Value = 376
$ws4.Range("F45").Value = 404
$ws4.Range("F46").Value = 836
$ws4.Range("F48").Value = 327
$ws4.Range("F49").Value = 6
